$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.080.55"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "3.429.26"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "409.17"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.12"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("E7").Value = "  +6.67%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  +7.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("E10").Value = "  +5.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.84"
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000224"
$ws.Range("E12").Value = "  +52.50%  "
$ws.Range("E13").Value = "  +10.62%  "
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "3.973.35"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.32"
$ws.Range("E16").Value = "  +7.59%  "
$ws.Range("D17").Value = "3.422.70"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.51"
$ws.Range("E18").Value = "  +7.40%  "
$ws.Range("E19").Value = "  +7.83%  "
$ws.Range("D20").Value = "62.019.52"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "456.54"
$ws.Range("E21").Value = "  +46.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "91.44"
$ws.Range("E22").Value = "  +8.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.23"
$ws.Range("E23").Value = "  +1.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.09"
$ws.Range("E24").Value = "  +2.53%  "
$ws.Range("E25").Value = "  +2.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "33.11"
$ws.Range("E26").Value = "  +11.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.18"
$ws.Range("E27").Value = "  +12.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.77"
$ws.Range("E28").Value = "  +1.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.70"
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("E30").Value = "  +6.50%  "
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.114"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.02"
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0503"
$ws.Range("E36").Value = "  +3.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.38"
$ws.Range("E37").Value = "  +5.04%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.38"
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("E40").Value = "  +7.81%  "
$ws.Range("E41").Value = "  -1.38%  "
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "142.56"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.28"
$ws.Range("E44").Value = "  +9.93%  "
$ws.Range("E45").Value = "  +1.30%  "
$ws.Range("E46").Value = "  +13.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.66"
$ws.Range("E47").Value = "  -0.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.33"
$ws.Range("E48").Value = "  +5.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.12"
$ws.Range("E49").Value = "  +8.49%  "
$ws.Range("E50").Value = "  +17.04%  "
$ws.Range("D51").Value = "3.777.82"
$ws.Range("E51").Value = "  -0.21%  "
